# Insert a new snapshot column right before the "nom" / "url_produit"
# columns (old BF/BG), shifting them one column to the right (BF->BG,
# BG->BH). This mirrors the periodic scraper appending a new timestamped
# price snapshot column to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at BF; everything from BF onward (nom, url_produit)
# shifts right by one (BF->BG, BG->BH). Excel also grows the used range /
# dimension automatically (A1:BG206 -> A1:BH206).
$ws.Columns("BF").Insert()

# Header row: the new column BF1 gets the new snapshot timestamp.
$ws.Range("BF1").Value2 = "2026-01-30 08:24:55"

# Data rows: the new BF column is a duplicate of the latest price already
# recorded in BE for that row (rows 2-80 have a numeric BE value). Rows
# 81-206 have no BE price recorded this round, so their new BF cell is
# left blank, matching BE.
for ($r = 2; $r -le 206; $r++) {
    $priceCell = $ws.Cells.Item($r, 57)   # column BE
    $newCell = $ws.Cells.Item($r, 58)     # column BF (newly inserted)
    $price = $priceCell.Value2
    if ($price -ne $null -and $price -ne "") {
        $newCell.Value2 = $price
    }
}
